$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: turn the "Vary_Total4000_Color4000" row into a
#     "Vary_Unlimited_Color1" row (matching the other Vary_* rows' shape) ---

# A2: name of the scenario
$ws.Range("A2").Value = "Vary_Unlimited_Color1"

# D2: total-quota select value (also needs the font/style used by the
# "definitionTotalSelectControl-listbox-item-1" cells elsewhere on the sheet,
# e.g. the style already present on A2/B2/F2/I2)
$ws.Range("D2").Value = "definitionTotalSelectControl-listbox-item-1"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null

# G2: total-quota value becomes the text "Unlimited" (was numeric 4000)
$ws.Range("G2").Value = "Unlimited"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null

# H2: numeric total quota count drops from 4000 to 1
$ws.Range("H2").Value = 1

# J2: color-quota value becomes the text "Unlimited" (was numeric 4000)
$ws.Range("J2").Value = "Unlimited"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null

# L2: numeric color quota count drops from 4000 to 1
$ws.Range("L2").Value = 1

$excel.CutCopyMode = $false

# --- Sheet2 no longer needs its own duplicate "Vary_Unlimited_Color1" row;
#     that data now lives in Sheet1!row 2. Remove row 12 from Sheet2 without
#     shifting the rows below it (delete then re-insert a blank row). ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows("12").Delete() | Out-Null
$ws2.Rows("15").Insert() | Out-Null

# --- sheet view: update the active selection ---
$ws.Range("G11").Select()
